$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the "Congrats, you are ready to go! ..." paragraph with the
#    following "Hello Pfizer!" paragraph into a single paragraph, inserting
#    proofErr markers and re-splitting the text into the runs seen in the
#    target XML.
# ---------------------------------------------------------------------------
$p1 = $null
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq "Congrats, you are ready to go! Your forked repository is setup and you're ready to develop.") {
        $p1 = $p
    }
}
$p2 = $p1.Next()

$mergeRange = $d.Range($p1.Range.Start, $p2.Range.End)
$mergeXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:proofErr w:type="gramStart"/>
<w:r><w:t xml:space="preserve">Congrats, </w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>Hello</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> Pfizer!</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>you are ready to go! Your forked repository is setup and you''re ready to develop.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>'
$mergeRange.InsertXML($mergeXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Table cards row: "[Placeholder Image: Person walking <break> towards
#    light]" -> merge into a single run, no page break.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("[Placeholder Image: Person walking towards light]", $true, $false, $false, $false, $false, $true, 1, $false, "[Placeholder Image: Person walking towards light]", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Table cards row: drop the lastRenderedPageBreak before "Uncertainty
#    eliminated".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Uncertainty eliminated", $true, $false, $false, $false, $false, $true, 1, $false, "Uncertainty eliminated", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Table cards row: re-flow "[Placeholder Text: Description of <break>
#    uncertainty eliminated]" -> "[Placeholder Text: Description of
#    uncertainty <break> eliminated]".
# ---------------------------------------------------------------------------
$target = $null
foreach ($t in $d.Tables) {
    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        for ($c = 1; $c -le $t.Columns.Count; $c++) {
            $cellText = $t.Cell($r, $c).Range.Text.TrimEnd([char]13, [char]7)
            if ($cellText -eq "[Placeholder Text: Description of uncertainty eliminated]") {
                $target = $t.Cell($r, $c).Range.Paragraphs(1)
            }
        }
    }
}

$descXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">[Placeholder Text: Description of uncertainty </w:t></w:r>
<w:r><w:lastRenderedPageBreak/><w:t>eliminated]</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>'
$target.Range.InsertXML($descXml) | Out-Null
